# Update cryptocurrency price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold plain numeric-looking text (e.g. "594.29"); force the
# NumberFormat to Text first so the COM .Value setter does not auto-coerce them
# into actual numbers (matching the source data which is text throughout column D).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.360.89"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "3.252.62"
$ws.Range("E3").Value = "  +3.92%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "594.29"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "140.20"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.246.84"
$ws.Range("E8").Value = "  +3.70%  "
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "3.794.04"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "3.249.82"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("D18").Value = "63.411.80"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "474.46"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("E22").Value = "  +3.55%  "
$ws.Range("D23").Value = "7.93"
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("D24").Value = "83.86"
$ws.Range("E24").Value = "  -4.31%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  +3.69%  "
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("D31").Value = "27.59"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("D34").Value = "2.52"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("D37").Value = "52.67"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").Value = "419.43"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("D41").Value = "2.983.76"
$ws.Range("E41").Value = "  +2.95%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("E43").Value = "  -5.48%  "
$ws.Range("E44").Value = "  -7.84%  "
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D48").Value = "25.76"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").Value = "122.03"
$ws.Range("E51").Value = "  +1.23%  "

$wb.Save()
